# Insert a new weekly data row for Brócoli - Terminal Hortofrutícola Agro Chillán
# at row 401, pushing the existing rows 401-430 down to 402-431.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 401 (shifts rows 401:430 down to 402:431)
$ws.Rows("401:401").Insert()

# Populate the newly inserted row 401 with the new weekly record
$ws.Range("A401").Value = 7
$ws.Range("B401").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C401").Value = "Ñuble"
$ws.Range("D401").Value = 45013
$ws.Range("E401").Value = 16
$ws.Range("F401").Value = 100112023
$ws.Range("G401").Value = "Brócoli"
$ws.Range("H401").Value = "Sin especificar"
$ws.Range("I401").Value = "Primera"
$ws.Range("J401").Value = 110
$ws.Range("K401").Value = 1000
$ws.Range("L401").Value = 1200
$ws.Range("M401").Value = 1109
$ws.Range("N401").Value = "$/unidad"
$ws.Range("O401").Value = "Región del Maule"
$ws.Range("P401").Value = 1109
$ws.Range("Q401").Value = 1
$ws.Range("R401").Value = "Hortaliza"
